$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1616.5938
$ws.Range("J17").Value = 1651.5172
$ws.Range("L17").Value = 4954.5516
$ws.Range("N17").Value = -5290.5516
$ws.Range("H81").Value = 80000
$ws.Range("J81").Value = 80000
$ws.Range("L81").Value = 80000
$ws.Range("N81").Value = -81996
$ws.Range("H84").Value = 80000
$ws.Range("J84").Value = 80000
$ws.Range("L84").Value = 240000
$ws.Range("N84").Value = -249984
$ws.Range("H100").Value = 69945.336
$ws.Range("I100").Value = 77451
$ws.Range("K100").Value = 77451
$ws.Range("M100").Value = -76910
$ws.Range("H112").Value = 7170.5093
$ws.Range("J112").Value = 7365.6226
$ws.Range("L112").Value = 22096.8678
$ws.Range("N112").Value = -24312.8678
$ws.Range("H131").Value = 2201.4167
$ws.Range("I131").Value = 1765.1818
$ws.Range("J131").Value = 7000
$ws.Range("K131").Value = 5295.5454
$ws.Range("L131").Value = 21000
$ws.Range("M131").Value = -255.5454
$ws.Range("N131").Value = -31080
$ws.Range("H137").Value = 11776.333
$ws.Range("I137").Value = 16621.637
$ws.Range("J137").Value = 4162.2856
$ws.Range("K137").Value = 49864.91099999999
$ws.Range("L137").Value = 12486.8568
$ws.Range("M137").Value = -47314.91099999999
$ws.Range("N137").Value = -17586.8568
$ws.Range("H138").Value = 2668.383
$ws.Range("I138").Value = 1836.6
$ws.Range("J138").Value = 5094.4165
$ws.Range("K138").Value = 5509.799999999999
$ws.Range("L138").Value = 15283.2495
$ws.Range("M138").Value = -369.7999999999993
$ws.Range("N138").Value = -25563.2495

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2700.8616
$ws.Range("I61").Value = 1833.2727
$ws.Range("K61").Value = 1833.2727
$ws.Range("M61").Value = -1621.2727
$ws.Range("H97").Value = 2723.425
$ws.Range("I97").Value = 2436.4333
$ws.Range("K97").Value = 2436.4333
$ws.Range("M97").Value = -1940.4333
$ws.Range("H110").Value = 1241.3462
$ws.Range("I110").Value = 1114.1818
$ws.Range("K110").Value = 1114.1818
$ws.Range("M110").Value = 930.8181999999999
$ws.Range("H122").Value = 1730.58
$ws.Range("I122").Value = 1650.375
$ws.Range("K122").Value = 4951.125
$ws.Range("M122").Value = -2501.125
$ws.Range("H132").Value = 1682.1428
$ws.Range("I132").Value = 1626.4688
$ws.Range("K132").Value = 4879.4064
$ws.Range("M132").Value = -2349.4064
$ws.Range("H136").Value = 2700.8616
$ws.Range("I136").Value = 1833.2727
$ws.Range("K136").Value = 5499.8181
$ws.Range("M136").Value = -2949.8181

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 77175
$ws.Range("I20").Value = 78918.62
$ws.Range("K20").Value = 78918.62
$ws.Range("M20").Value = -78671.62
$ws.Range("H50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("N50").ClearContents()
$ws.Range("H105").Value = 1356.0834
$ws.Range("I105").Value = 1307.1428
$ws.Range("K105").Value = 1307.1428
$ws.Range("M105").Value = 439.8571999999999
$ws.Range("H134").Value = 2661.85
$ws.Range("I134").Value = 2508.0408
$ws.Range("K134").Value = 7524.1224
$ws.Range("M134").Value = -4989.1224

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 624.4286
$ws.Range("I22").Value = 550.5
$ws.Range("J22").Value = 723
$ws.Range("K22").Value = 550.5
$ws.Range("L22").Value = 723
$ws.Range("M22").Value = -200.5
$ws.Range("N22").Value = -1423
$ws.Range("H59").Value = 42447.4
$ws.Range("J59").Value = 42447.4
$ws.Range("L59").Value = 42447.4
$ws.Range("N59").Value = -44737.4
$ws.Range("H87").Value = 35000
$ws.Range("J87").Value = 35000
$ws.Range("L87").Value = 35000
$ws.Range("N87").Value = -37372
$ws.Range("H90").Value = 35000
$ws.Range("J90").Value = 35000
$ws.Range("L90").Value = 105000
$ws.Range("N90").Value = -116856
$ws.Range("H105").Value = 1745.6923
$ws.Range("I105").Value = 1824.6666
$ws.Range("J105").Value = 798
$ws.Range("K105").Value = 1824.6666
$ws.Range("L105").Value = 798
$ws.Range("M105").Value = -77.66660000000002
$ws.Range("N105").Value = -4292
$ws.Range("H132").Value = 70149.336
$ws.Range("I132").Value = 103055.664
$ws.Range("J132").Value = 4336.6665
$ws.Range("K132").Value = 309166.992
$ws.Range("L132").Value = 13009.9995
$ws.Range("M132").Value = -306636.992
$ws.Range("N132").Value = -18069.9995

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 44618320
$ws.Range("I4").Value = 52716788
$ws.Range("K4").Value = 158150364
$ws.Range("M4").Value = -158150252
$ws.Range("H6").Value = 4.285714
$ws.Range("I6").Value = 4.6666665
$ws.Range("K6").Value = 13.9999995
$ws.Range("M6").Value = 99.0000005
$ws.Range("H7").Value = 283
$ws.Range("I7").Value = 250
$ws.Range("K7").Value = 750
$ws.Range("M7").Value = -638
$ws.Range("H107").Value = 989.7692
$ws.Range("I107").Value = 1197.4286
$ws.Range("J107").Value = 747.5
$ws.Range("K107").Value = 3592.2858
$ws.Range("L107").Value = 2242.5
$ws.Range("M107").Value = -1672.2858
$ws.Range("N107").Value = -6082.5
$ws.Range("H129").Value = 3419
$ws.Range("I129").Value = 1590
$ws.Range("J129").Value = 4594.7856
$ws.Range("K129").Value = 4770
$ws.Range("L129").Value = 13784.3568
$ws.Range("M129").Value = 230
$ws.Range("N129").Value = -23784.3568
$ws.Range("H131").Value = 63455.973
$ws.Range("J131").Value = 2182.7234
$ws.Range("L131").Value = 6548.1702
$ws.Range("N131").Value = -16628.1702
$ws.Range("H132").Value = 1790.6765
$ws.Range("I132").Value = 1851.4
$ws.Range("J132").Value = 1742.7368
$ws.Range("K132").Value = 16662.6
$ws.Range("L132").Value = 15684.6312
$ws.Range("M132").Value = -14132.6
$ws.Range("N132").Value = -20744.6312
$ws.Range("H137").Value = 2193.8572
$ws.Range("I137").Value = 2193.8572
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 6581.571599999999
$ws.Range("L137").Value = 0
$ws.Range("M137").Value = -1481.571599999999
$ws.Range("N137").ClearContents()

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5774.4443
$ws.Range("I70").Value = 3994
$ws.Range("K70").Value = 3994
$ws.Range("M70").Value = -3724
$ws.Range("H73").Value = 5774.4443
$ws.Range("I73").Value = 3994
$ws.Range("K73").Value = 3994
$ws.Range("M73").Value = -3058

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 4802.3335
$ws.Range("I46").Value = 2050.125
$ws.Range("K46").Value = 2050.125
$ws.Range("M46").Value = -1862.125
$ws.Range("H68").Value = 3295.3333
$ws.Range("I68").Value = 2754.4
$ws.Range("J68").Value = 6000
$ws.Range("K68").Value = 2754.4
$ws.Range("L68").Value = 6000
$ws.Range("M68").Value = -2005.4
$ws.Range("N68").Value = -7498
$ws.Range("H71").Value = 3295.3333
$ws.Range("I71").Value = 2754.4
$ws.Range("J71").Value = 6000
$ws.Range("K71").Value = 13772
$ws.Range("L71").Value = 30000
$ws.Range("M71").Value = -10028
$ws.Range("N71").Value = -37488
$ws.Range("H93").Value = 2926.5
$ws.Range("I93").Value = 2332
$ws.Range("J93").Value = 4264.125
$ws.Range("K93").Value = 2332
$ws.Range("L93").Value = 4264.125
$ws.Range("M93").Value = -1084
$ws.Range("N93").Value = -6760.125
$ws.Range("H122").Value = 29972.637
$ws.Range("I122").Value = 33054.11
$ws.Range("J122").Value = 16106
$ws.Range("K122").Value = 99162.33
$ws.Range("L122").Value = 48318
$ws.Range("M122").Value = -96712.33
$ws.Range("N122").Value = -53218
$ws.Range("H141").Value = 110700
$ws.Range("J141").Value = 110700
$ws.Range("L141").Value = 110700
$ws.Range("N141").Value = -121060

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1039.7
$ws.Range("I96").Value = 1316.3334
$ws.Range("J96").Value = 624.75
$ws.Range("K96").Value = 1316.3334
$ws.Range("L96").Value = 624.75
$ws.Range("M96").Value = 56.66660000000002
$ws.Range("N96").Value = -3370.75
$ws.Range("H132").Value = 4251951.5
$ws.Range("I132").Value = 4325175
$ws.Range("K132").Value = 12975525
$ws.Range("M132").Value = -12972995
$ws.Range("H136").Value = 16357.355
$ws.Range("I136").Value = 16722.295
$ws.Range("J136").Value = 300
$ws.Range("K136").Value = 50166.88499999999
$ws.Range("L136").Value = 900
$ws.Range("M136").Value = -47616.88499999999
$ws.Range("N136").Value = -6000
$ws.Range("H141").Value = 102115.71
$ws.Range("J141").Value = 102115.71
$ws.Range("L141").Value = 102115.71
$ws.Range("N141").Value = -112475.71
